$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 422.1640277809471
$ws.Range("B2").Value = 72.29138580545695

$ws.Range("A3").Value = 434.0774227789693
$ws.Range("B3").Value = 68.84493981932113

$ws.Range("A4").Value = 438.4150340832591
$ws.Range("B4").Value = 70.97626561410722

$ws.Range("A5").Value = 431.0540579015487
$ws.Range("B5").Value = 68.03717310698939

$ws.Range("A6").Value = 438.7164927469763
$ws.Range("B6").Value = 54.37358433650676

$ws.Range("A7").Value = 435.4090084474875
$ws.Range("B7").Value = 79.80926178634144

$ws.Range("A8").Value = 421.0657516171254
$ws.Range("B8").Value = 69.26105417588337

$ws.Range("A9").Value = 409.926635397423
$ws.Range("B9").Value = 56.87163956602839

$ws.Range("A10").Value = 439.215700756011
$ws.Range("B10").Value = 57.72587477486562

$ws.Range("A11").Value = 432.7363602554482
$ws.Range("B11").Value = 74.92266025006445
